$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 26 ("S18" step) following the same pattern/formatting as row 25.
# Copy formatting from the last populated row so the new row's styles match.
$ws.Range("A25:F25").Copy()
$ws.Range("A26:F26").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A26").Value = "Final feature extraction consensus sheet"
$ws.Range("B26").Value = 15.0
$ws.Range("C26").Value = "11/6/2023"
$ws.Range("D26").Value = 0.0
$ws.Range("E26").Value = 73.0
$ws.Range("F26").Value = "S18"
